# Insert a new data row at row 71 (pushing the existing rows 71-148 down to
# 72-149) and populate it with the new "Femacal de La Calera - Haba" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(71).Insert()

$ws.Range("A71").Value2 = 3
$ws.Range("B71").Value2 = "Femacal de La Calera"
$ws.Range("C71").Value2 = "Coquimbo"
$ws.Range("D71").Value2 = 44741
$ws.Range("E71").Value2 = 5
$ws.Range("F71").Value2 = 100112026
$ws.Range("G71").Value2 = "Haba"
$ws.Range("H71").Value2 = "Sin especificar"
$ws.Range("I71").Value2 = "Primera"
$ws.Range("J71").Value2 = 73
$ws.Range("K71").Value2 = 21000
$ws.Range("L71").Value2 = 22000
$ws.Range("M71").Value2 = 21479
$ws.Range("N71").Value2 = '$/saco 25 kilos'
$ws.Range("O71").Value2 = "Provincia de Limarí"
$ws.Range("P71").Value2 = 859
$ws.Range("Q71").Value2 = 25
$ws.Range("R71").Value2 = "Hortaliza"
